$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# The workbook currently has 3 sheets: 2021-Q3, 2021-Q4, 总计.
# The edit inserts a new quarter (2022-Q1) of fund-holdings detail
# right before the running '总计' (grand total) summary sheet, and
# records its row in that summary.
#
# To land on the same sheetId/rId numbering as the target file
# (2022-Q1 -> sheetId 3/rId3, 总计 -> sheetId 4/rId4) we rename the
# existing '总计' sheet (sheetId 3) to '2022-Q1' and then add a
# brand new '总计' sheet right after it.
# ------------------------------------------------------------------
$q1 = $wb.Worksheets.Item(3)
$q1.Name = "2022-Q1"

# ClearContents (not Clear!) drops the old 3-row summary's values
# but keeps the cell formatting, so the header row (B1:D1) and
# column A (A2:A3) keep exactly the style index they already had;
# that style gets reused/extended onto the rest of the rebuilt
# sheet below instead of having to invent it from scratch.
$q1.Cells.ClearContents()

$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

# Match the sheetPr/pageMargins cosmetics the other sheets in this
# workbook already use (outline summary-below/right, 0.75in side
# margins, 1in top/bottom, 0.5in header/footer) instead of the
# engine's blank-sheet defaults.
$total.Outline.SummaryRow = 1
$total.Outline.SummaryColumn = 1
$total.PageSetup.LeftMargin = 54
$total.PageSetup.RightMargin = 54
$total.PageSetup.TopMargin = 72
$total.PageSetup.BottomMargin = 72
$total.PageSetup.HeaderMargin = 36
$total.PageSetup.FooterMargin = 36

# ------------------------------------------------------------------
# 2022-Q1 sheet: header row (B1:H1) + 27 fund-holding detail rows
# (A2:H28). Column A holds a 0-based row index; B..G are text-typed
# fund codes/names/numbers (matching the source data's inlineStr
# cells); H is a plain integer rank.
# ------------------------------------------------------------------
$q1.Cells.Item(1,2).Value = "基金代码"
$q1.Cells.Item(1,3).Value = "基金名称"
$q1.Cells.Item(1,4).Value = "基金规模"
$q1.Cells.Item(1,5).Value = "股票总仓位"
$q1.Cells.Item(1,6).Value = "仓位占比"
$q1.Cells.Item(1,7).Value = "持有市值(亿元)"
$q1.Cells.Item(1,8).Value = "仓位排名"

# Stamp the existing header style (still on B1/C1/D1 after
# ClearContents) onto the newly used header cells E1:H1.
$q1.Cells.Item(1,2).Copy()
$q1.Range("E1:H1").PasteSpecial(-4122)

# Stamp the existing column-A style (still on A2/A3 after
# ClearContents) onto the rest of column A (A4:A28).
$q1.Cells.Item(2,1).Copy()
$q1.Range("A4:A28").PasteSpecial(-4122)

# Force columns B:G to store plain text (matches the source file's
# t="inlineStr" cells) instead of having numeric-looking values
# (fund codes like 010379, decimals like 55.11) auto-coerced to
# numbers / losing leading zeros. The style gets dropped back to
# the default afterwards so no stray style index is left behind
# on these data cells (only column A is deliberately styled).
$textRange = $q1.Range("B2:G28")
$textRange.NumberFormat = "@"

$q1.Cells.Item(2,1).Value = 0
$q1.Cells.Item(2,2).Value = "010379"
$q1.Cells.Item(2,3).Value = "广发均衡优选混合A"
$q1.Cells.Item(2,4).Value = "55.11"
$q1.Cells.Item(2,5).Value = "64.69"
$q1.Cells.Item(2,6).Value = "4.61"
$q1.Cells.Item(2,7).Value = "2.5406"
$q1.Cells.Item(2,8).Value = 8

$q1.Cells.Item(3,1).Value = 1
$q1.Cells.Item(3,2).Value = "008297"
$q1.Cells.Item(3,3).Value = "广发价值优势混合"
$q1.Cells.Item(3,4).Value = "28.25"
$q1.Cells.Item(3,5).Value = "93.96"
$q1.Cells.Item(3,6).Value = "5.52"
$q1.Cells.Item(3,7).Value = "1.5594"
$q1.Cells.Item(3,8).Value = 5

$q1.Cells.Item(4,1).Value = 2
$q1.Cells.Item(4,2).Value = "009887"
$q1.Cells.Item(4,3).Value = "广发稳健优选六个月持有期混合A"
$q1.Cells.Item(4,4).Value = "30.91"
$q1.Cells.Item(4,5).Value = "64.69"
$q1.Cells.Item(4,6).Value = "4.65"
$q1.Cells.Item(4,7).Value = "1.4373"
$q1.Cells.Item(4,8).Value = 6

$q1.Cells.Item(5,1).Value = 3
$q1.Cells.Item(5,2).Value = "014591"
$q1.Cells.Item(5,3).Value = "广发瑞誉一年持有期混合A"
$q1.Cells.Item(5,4).Value = "38.60"
$q1.Cells.Item(5,5).Value = "93.40"
$q1.Cells.Item(5,6).Value = "3.71"
$q1.Cells.Item(5,7).Value = "1.4321"
$q1.Cells.Item(5,8).Value = 9

$q1.Cells.Item(6,1).Value = 4
$q1.Cells.Item(6,2).Value = "011194"
$q1.Cells.Item(6,3).Value = "广发睿铭两年持有期混合型证券投资基金A"
$q1.Cells.Item(6,4).Value = "20.55"
$q1.Cells.Item(6,5).Value = "74.56"
$q1.Cells.Item(6,6).Value = "5.40"
$q1.Cells.Item(6,7).Value = "1.1097"
$q1.Cells.Item(6,8).Value = 5

$q1.Cells.Item(7,1).Value = 5
$q1.Cells.Item(7,2).Value = "398001"
$q1.Cells.Item(7,3).Value = "中海优质成长混合"
$q1.Cells.Item(7,4).Value = "14.42"
$q1.Cells.Item(7,5).Value = "90.86"
$q1.Cells.Item(7,6).Value = "6.63"
$q1.Cells.Item(7,7).Value = "0.9560"
$q1.Cells.Item(7,8).Value = 2

$q1.Cells.Item(8,1).Value = 6
$q1.Cells.Item(8,2).Value = "270022"
$q1.Cells.Item(8,3).Value = "广发内需增长混合A"
$q1.Cells.Item(8,4).Value = "15.92"
$q1.Cells.Item(8,5).Value = "79.56"
$q1.Cells.Item(8,6).Value = "5.78"
$q1.Cells.Item(8,7).Value = "0.9202"
$q1.Cells.Item(8,8).Value = 5

$q1.Cells.Item(9,1).Value = 7
$q1.Cells.Item(9,2).Value = "009888"
$q1.Cells.Item(9,3).Value = "广发稳健优选六个月持有期混合C"
$q1.Cells.Item(9,4).Value = "11.86"
$q1.Cells.Item(9,5).Value = "64.69"
$q1.Cells.Item(9,6).Value = "4.65"
$q1.Cells.Item(9,7).Value = "0.5515"
$q1.Cells.Item(9,8).Value = 6

$q1.Cells.Item(10,1).Value = 8
$q1.Cells.Item(10,2).Value = "011134"
$q1.Cells.Item(10,3).Value = "广发价值优选混合A"
$q1.Cells.Item(10,4).Value = "6.35"
$q1.Cells.Item(10,5).Value = "93.95"
$q1.Cells.Item(10,6).Value = "6.52"
$q1.Cells.Item(10,7).Value = "0.4140"
$q1.Cells.Item(10,8).Value = 5

$q1.Cells.Item(11,1).Value = 9
$q1.Cells.Item(11,2).Value = "011195"
$q1.Cells.Item(11,3).Value = "广发睿铭两年持有期混合型证券投资基金C"
$q1.Cells.Item(11,4).Value = "7.63"
$q1.Cells.Item(11,5).Value = "74.56"
$q1.Cells.Item(11,6).Value = "5.40"
$q1.Cells.Item(11,7).Value = "0.4120"
$q1.Cells.Item(11,8).Value = 5

$q1.Cells.Item(12,1).Value = 10
$q1.Cells.Item(12,2).Value = "002430"
$q1.Cells.Item(12,3).Value = "中银丰利灵活配置混合A"
$q1.Cells.Item(12,4).Value = "8.76"
$q1.Cells.Item(12,5).Value = "28.78"
$q1.Cells.Item(12,6).Value = "1.85"
$q1.Cells.Item(12,7).Value = "0.1621"
$q1.Cells.Item(12,8).Value = 10

$q1.Cells.Item(13,1).Value = 11
$q1.Cells.Item(13,2).Value = "010380"
$q1.Cells.Item(13,3).Value = "广发均衡优选混合C"
$q1.Cells.Item(13,4).Value = "3.49"
$q1.Cells.Item(13,5).Value = "64.69"
$q1.Cells.Item(13,6).Value = "4.61"
$q1.Cells.Item(13,7).Value = "0.1609"
$q1.Cells.Item(13,8).Value = 8

$q1.Cells.Item(14,1).Value = 12
$q1.Cells.Item(14,2).Value = "014592"
$q1.Cells.Item(14,3).Value = "广发瑞誉一年持有期混合C"
$q1.Cells.Item(14,4).Value = "4.08"
$q1.Cells.Item(14,5).Value = "93.40"
$q1.Cells.Item(14,6).Value = "3.71"
$q1.Cells.Item(14,7).Value = "0.1514"
$q1.Cells.Item(14,8).Value = 9

$q1.Cells.Item(15,1).Value = 13
$q1.Cells.Item(15,2).Value = "002616"
$q1.Cells.Item(15,3).Value = "中银益利灵活配置混合A"
$q1.Cells.Item(15,4).Value = "5.90"
$q1.Cells.Item(15,5).Value = "29.85"
$q1.Cells.Item(15,6).Value = "2.17"
$q1.Cells.Item(15,7).Value = "0.1280"
$q1.Cells.Item(15,8).Value = 6

$q1.Cells.Item(16,1).Value = 14
$q1.Cells.Item(16,2).Value = "003850"
$q1.Cells.Item(16,3).Value = "中银锦利灵活配置混合A"
$q1.Cells.Item(16,4).Value = "6.37"
$q1.Cells.Item(16,5).Value = "28.68"
$q1.Cells.Item(16,6).Value = "1.75"
$q1.Cells.Item(16,7).Value = "0.1115"
$q1.Cells.Item(16,8).Value = 10

$q1.Cells.Item(17,1).Value = 15
$q1.Cells.Item(17,2).Value = "011135"
$q1.Cells.Item(17,3).Value = "广发价值优选混合C"
$q1.Cells.Item(17,4).Value = "1.48"
$q1.Cells.Item(17,5).Value = "93.95"
$q1.Cells.Item(17,6).Value = "6.52"
$q1.Cells.Item(17,7).Value = "0.0965"
$q1.Cells.Item(17,8).Value = 5

$q1.Cells.Item(18,1).Value = 16
$q1.Cells.Item(18,2).Value = "002431"
$q1.Cells.Item(18,3).Value = "中银丰利灵活配置混合C"
$q1.Cells.Item(18,4).Value = "2.24"
$q1.Cells.Item(18,5).Value = "28.78"
$q1.Cells.Item(18,6).Value = "1.85"
$q1.Cells.Item(18,7).Value = "0.0414"
$q1.Cells.Item(18,8).Value = 10

$q1.Cells.Item(19,1).Value = 17
$q1.Cells.Item(19,2).Value = "010740"
$q1.Cells.Item(19,3).Value = "汇安核心价值混合A"
$q1.Cells.Item(19,4).Value = "0.95"
$q1.Cells.Item(19,5).Value = "93.68"
$q1.Cells.Item(19,6).Value = "3.07"
$q1.Cells.Item(19,7).Value = "0.0292"
$q1.Cells.Item(19,8).Value = 10

$q1.Cells.Item(20,1).Value = 18
$q1.Cells.Item(20,2).Value = "003851"
$q1.Cells.Item(20,3).Value = "中银锦利灵活配置混合C"
$q1.Cells.Item(20,4).Value = "1.25"
$q1.Cells.Item(20,5).Value = "28.68"
$q1.Cells.Item(20,6).Value = "1.75"
$q1.Cells.Item(20,7).Value = "0.0219"
$q1.Cells.Item(20,8).Value = 10

$q1.Cells.Item(21,1).Value = 19
$q1.Cells.Item(21,2).Value = "001252"
$q1.Cells.Item(21,3).Value = "中海进取收益灵活配置混合"
$q1.Cells.Item(21,4).Value = "0.23"
$q1.Cells.Item(21,5).Value = "92.65"
$q1.Cells.Item(21,6).Value = "6.20"
$q1.Cells.Item(21,7).Value = "0.0143"
$q1.Cells.Item(21,8).Value = 1

$q1.Cells.Item(22,1).Value = 20
$q1.Cells.Item(22,2).Value = "002617"
$q1.Cells.Item(22,3).Value = "中银益利灵活配置混合C"
$q1.Cells.Item(22,4).Value = "0.58"
$q1.Cells.Item(22,5).Value = "29.85"
$q1.Cells.Item(22,6).Value = "2.17"
$q1.Cells.Item(22,7).Value = "0.0126"
$q1.Cells.Item(22,8).Value = 6

$q1.Cells.Item(23,1).Value = 21
$q1.Cells.Item(23,2).Value = "011677"
$q1.Cells.Item(23,3).Value = "中银睿丰回报混合型证券投资基金A"
$q1.Cells.Item(23,4).Value = "0.73"
$q1.Cells.Item(23,5).Value = "20.29"
$q1.Cells.Item(23,6).Value = "1.67"
$q1.Cells.Item(23,7).Value = "0.0122"
$q1.Cells.Item(23,8).Value = 6

$q1.Cells.Item(24,1).Value = 22
$q1.Cells.Item(24,2).Value = "010741"
$q1.Cells.Item(24,3).Value = "汇安核心价值混合C"
$q1.Cells.Item(24,4).Value = "0.36"
$q1.Cells.Item(24,5).Value = "93.68"
$q1.Cells.Item(24,6).Value = "3.07"
$q1.Cells.Item(24,7).Value = "0.0111"
$q1.Cells.Item(24,8).Value = 10

$q1.Cells.Item(25,1).Value = 23
$q1.Cells.Item(25,2).Value = "011183"
$q1.Cells.Item(25,3).Value = "广发内需增长混合C"
$q1.Cells.Item(25,4).Value = "0.12"
$q1.Cells.Item(25,5).Value = "79.56"
$q1.Cells.Item(25,6).Value = "5.78"
$q1.Cells.Item(25,7).Value = "0.0069"
$q1.Cells.Item(25,8).Value = 5

$q1.Cells.Item(26,1).Value = 24
$q1.Cells.Item(26,2).Value = "000822"
$q1.Cells.Item(26,3).Value = "东海美丽中国灵活配置混合"
$q1.Cells.Item(26,4).Value = "0.04"
$q1.Cells.Item(26,5).Value = "89.82"
$q1.Cells.Item(26,6).Value = "2.07"
$q1.Cells.Item(26,7).Value = "0.0008"
$q1.Cells.Item(26,8).Value = 9

$q1.Cells.Item(27,1).Value = 25
$q1.Cells.Item(27,2).Value = "001797"
$q1.Cells.Item(27,3).Value = "华融新利灵活配置混合"
$q1.Cells.Item(27,4).Value = "0.02"
$q1.Cells.Item(27,5).Value = "48.66"
$q1.Cells.Item(27,6).Value = "2.91"
$q1.Cells.Item(27,7).Value = "0.0006"
$q1.Cells.Item(27,8).Value = 3

$q1.Cells.Item(28,1).Value = 26
$q1.Cells.Item(28,2).Value = "011678"
$q1.Cells.Item(28,3).Value = "中银睿丰回报混合型证券投资基金C"
$q1.Cells.Item(28,4).Value = "0.00"
$q1.Cells.Item(28,5).Value = "20.29"
$q1.Cells.Item(28,6).Value = "1.67"
$q1.Cells.Item(28,8).Value = 6

# Reset B:G back to the default style (keeps the text typing that
# was just written, drops the "@" number-format style index).
$textRange.Style = "Normal"

# G28 is a genuine number (0), not text -- set it now that
# the range is back to a plain/default (non-text) format.
$q1.Cells.Item(28,7).Value = 0

# ------------------------------------------------------------------
# 总计 (grand total) sheet: same header as before, with a new row
# for 2022-Q1 inserted above the existing 2021-Q4 / 2021-Q3 rows.
# ------------------------------------------------------------------
$total.Cells.Item(1,2).Value = "日期"
$total.Cells.Item(1,3).Value = "持有数量(只)"
$total.Cells.Item(1,4).Value = "持有市值(亿元)"

$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q1"
$total.Cells.Item(2,3).Value = 27
$total.Cells.Item(2,4).Value = 12.29

$total.Cells.Item(3,1).Value = 1
$total.Cells.Item(3,2).Value = "2021-Q4"
$total.Cells.Item(3,3).Value = 20
$total.Cells.Item(3,4).Value = 17.36

$total.Cells.Item(4,1).Value = 2
$total.Cells.Item(4,2).Value = "2021-Q3"
$total.Cells.Item(4,3).Value = 14
$total.Cells.Item(4,4).Value = 1.64

# Style: header row (B1:D1) and column A (A2:A4) use the sheet's
# default 'index/header' style -- this is a brand new sheet so
# there is no pre-existing style-2 cell on it to copy from. Copy
# it over from the 2022-Q1 sheet, which already carries that
# exact style (on its own header/column-A cells).
$q1.Cells.Item(1,2).Copy()
$total.Range("B1:D1").PasteSpecial(-4122)
$q1.Cells.Item(2,1).Copy()
$total.Range("A2:A4").PasteSpecial(-4122)

# Re-activate the original first/selected tab -- adding a sheet
# above switched the active sheet to it, which would otherwise
# incorrectly move tabSelected/activeTab onto "总计".
$wb.Worksheets.Item(1).Activate()
